$d = $word.ActiveDocument

# The last paragraph in the document is the second "Eu,teste query..." block.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# 1) Insert a new, empty paragraph right after it (default paragraph properties).
$lastPara.Range.InsertParagraphAfter()

# 2) Insert another new paragraph after that empty one, carrying the same
#    "space after" formatting (5pt == w:after="100") as the paragraph above,
#    with the text "Meu teste de ato inicial" followed by a line break.
$emptyPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$emptyPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Format.SpaceAfter = 5
$newPara.Range.Text = "Meu teste de ato inicial"

$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertAfter([char]11)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ("[$i] '" + $p.Range.Text + "'")
}
